{"js": "// Replace the arithmetic-expression text in each cell of the single table in\n// the document body, cell-by-cell (row-major order), while preserving the\n// existing run/paragraph formatting (font, size, alignment).\n//\n// A single document-wide find/replace is not safe here: a handful of\n// expressions (e.g. \"90-53=\", \"33-5=\") occur in more than one cell but map to\n// *different* replacement text depending on which cell they're in, so every\n// cell is addressed individually by its zero-based (row, col) position, the\n// existing text is confirmed with a cell-scoped search, and only the matched\n// range's text is swapped \u2014 leaving the run's formatting untouched.\nconst cellEdits = [\n  { r: 0, c: 0, o: \"66+18=\", n: \"26+46=\" },\n  { r: 0, c: 1, o: \"48+38=\", n: \"39+17=\" },\n  { r: 0, c: 2, o: \"9+48=\", n: \"85+9=\" },\n  { r: 0, c: 3, o: \"3+78=\", n: \"19+19=\" },\n  { r: 0, c: 4, o: \"36-18=\", n: \"3+9=\" },\n  { r: 1, c: 0, o: \"24+58=\", n: \"50-12=\" },\n  { r: 1, c: 1, o: \"27+27=\", n: \"97-19=\" },\n  { r: 1, c: 2, o: \"82-74=\", n: \"55-6=\" },\n  { r: 1, c: 3, o: \"40-35=\", n: \"61-39=\" },\n  { r: 1, c: 4, o: \"49+44=\", n: \"29+47=\" },\n  { r: 2, c: 0, o: \"41-17=\", n: \"66+16=\" },\n  { r: 2, c: 1, o: \"60-59=\", n: \"92-76=\" },\n  { r: 2, c: 2, o: \"86-67=\", n: \"38+44=\" },\n  { r: 2, c: 3, o: \"14+18=\", n: \"73-29=\" },\n  { r: 2, c: 4, o: \"9+37=\", n: \"91-58=\" },\n  { r: 3, c: 0, o: \"5+7=\", n: \"40-22=\" },\n  { r: 3, c: 1, o: \"55-47=\", n: \"93-78=\" },\n  { r: 3, c: 2, o: \"25+38=\", n: \"85-37=\" },\n  { r: 3, c: 3, o: \"39+43=\", n: \"59+4=\" },\n  { r: 3, c: 4, o: \"66+19=\", n: \"9+62=\" },\n  { r: 4, c: 0, o: \"22-15=\", n: \"71-42=\" },\n  { r: 4, c: 1, o: \"92-18=\", n: \"57+24=\" },\n  { r: 4, c: 2, o: \"18+57=\", n: \"19+79=\" },\n  { r: 4, c: 3, o: \"80-29=\", n: \"19+25=\" },\n  { r: 4, c: 4, o: \"50-28=\", n: \"4+39=\" },\n  { r: 5, c: 0, o: \"77-38=\", n: \"78+13=\" },\n  { r: 5, c: 1, o: \"7+77=\", n: \"48+43=\" },\n  { r: 5, c: 2, o: \"9+64=\", n: \"45+39=\" },\n  { r: 5, c: 3, o: \"50-48=\", n: \"17+18=\" },\n  { r: 5, c: 4, o: \"49+35=\", n: \"16+75=\" },\n  { r: 6, c: 0, o: \"29+45=\", n: \"61-53=\" },\n  { r: 6, c: 1, o: \"9+18=\", n: \"63-24=\" },\n  { r: 6, c: 2, o: \"9+58=\", n: \"19+27=\" },\n  { r: 6, c: 3, o: \"17+55=\", n: \"43-36=\" },\n  { r: 6, c: 4, o: \"90-53=\", n: \"81-2=\" },\n  { r: 7, c: 0, o: \"19+66=\", n: \"36+7=\" },\n  { r: 7, c: 1, o: \"56+37=\", n: \"47+37=\" },\n  { r: 7, c: 2, o: \"73-36=\", n: \"26+6=\" },\n  { r: 7, c: 3, o: \"22-9=\", n: \"77+8=\" },\n  { r: 7, c: 4, o: \"60-8=\", n: \"26+29=\" },\n  { r: 8, c: 0, o: \"24+67=\", n: \"94-76=\" },\n  { r: 8, c: 1, o: \"33-5=\", n: \"9+76=\" },\n  { r: 8, c: 2, o: \"18+38=\", n: \"72-65=\" },\n  { r: 8, c: 3, o: \"86-8=\", n: \"33-29=\" },\n  { r: 8, c: 4, o: \"85-67=\", n: \"83-5=\" },\n  { r: 9, c: 0, o: \"21-7=\", n: \"84-59=\" },\n  { r: 9, c: 1, o: \"84-65=\", n: \"37-29=\" },\n  { r: 9, c: 2, o: \"89+9=\", n: \"73-38=\" },\n  { r: 9, c: 3, o: \"18+65=\", n: \"95-78=\" },\n  { r: 9, c: 4, o: \"94-86=\", n: \"74-48=\" },\n  { r: 10, c: 0, o: \"95-79=\", n: \"46+37=\" },\n  { r: 10, c: 1, o: \"85-36=\", n: \"90-32=\" },\n  { r: 10, c: 2, o: \"29+26=\", n: \"74-25=\" },\n  { r: 10, c: 3, o: \"62-38=\", n: \"6+49=\" },\n  { r: 10, c: 4, o: \"80-7=\", n: \"48+23=\" },\n  { r: 11, c: 0, o: \"8+13=\", n: \"61-6=\" },\n  { r: 11, c: 1, o: \"46+9=\", n: \"69+4=\" },\n  { r: 11, c: 2, o: \"39+42=\", n: \"61-44=\" },\n  { r: 11, c: 3, o: \"90-14=\", n: \"14-8=\" },\n  { r: 11, c: 4, o: \"29+35=\", n: \"82-17=\" },\n  { r: 12, c: 0, o: \"81-6=\", n: \"30-26=\" },\n  { r: 12, c: 1, o: \"73-9=\", n: \"9+68=\" },\n  { r: 12, c: 2, o: \"9+5=\", n: \"62-16=\" },\n  { r: 12, c: 3, o: \"17+49=\", n: \"9+57=\" },\n  { r: 12, c: 4, o: \"5+48=\", n: \"25+47=\" },\n  { r: 13, c: 0, o: \"92-3=\", n: \"19+55=\" },\n  { r: 13, c: 1, o: \"7+7=\", n: \"58+19=\" },\n  { r: 13, c: 2, o: \"16+25=\", n: \"70-41=\" },\n  { r: 13, c: 3, o: \"15+79=\", n: \"86-58=\" },\n  { r: 13, c: 4, o: \"75-7=\", n: \"59+22=\" },\n  { r: 14, c: 0, o: \"79+13=\", n: \"81-27=\" },\n  { r: 14, c: 1, o: \"90-53=\", n: \"36+28=\" },\n  { r: 14, c: 2, o: \"65+16=\", n: \"45+37=\" },\n  { r: 14, c: 3, o: \"33-6=\", n: \"47+6=\" },\n  { r: 14, c: 4, o: \"64+9=\", n: \"90-66=\" },\n  { r: 15, c: 0, o: \"41-37=\", n: \"33-9=\" },\n  { r: 15, c: 1, o: \"26+36=\", n: \"41-3=\" },\n  { r: 15, c: 2, o: \"7+87=\", n: \"57+28=\" },\n  { r: 15, c: 3, o: \"71-17=\", n: \"51-22=\" },\n  { r: 15, c: 4, o: \"46-28=\", n: \"16+29=\" },\n  { r: 16, c: 0, o: \"75-48=\", n: \"68+28=\" },\n  { r: 16, c: 1, o: \"36-7=\", n: \"23+8=\" },\n  { r: 16, c: 2, o: \"18+26=\", n: \"88+5=\" },\n  { r: 16, c: 3, o: \"98-29=\", n: \"80-2=\" },\n  { r: 16, c: 4, o: \"33-5=\", n: \"11-2=\" },\n  { r: 17, c: 0, o: \"27+24=\", n: \"64-47=\" },\n  { r: 17, c: 1, o: \"51-36=\", n: \"16+27=\" },\n  { r: 17, c: 2, o: \"55-27=\", n: \"70-61=\" },\n  { r: 17, c: 3, o: \"66+17=\", n: \"84-16=\" },\n  { r: 17, c: 4, o: \"35+47=\", n: \"87+7=\" },\n  { r: 18, c: 0, o: \"4+87=\", n: \"51-44=\" },\n  { r: 18, c: 1, o: \"67-29=\", n: \"38+18=\" },\n  { r: 18, c: 2, o: \"7+39=\", n: \"40-11=\" },\n  { r: 18, c: 3, o: \"35+26=\", n: \"61-22=\" },\n  { r: 18, c: 4, o: \"80-57=\", n: \"71-59=\" },\n  { r: 19, c: 0, o: \"63-56=\", n: \"18+19=\" },\n  { r: 19, c: 1, o: \"36+25=\", n: \"37+5=\" },\n  { r: 19, c: 2, o: \"3+28=\", n: \"30-7=\" },\n  { r: 19, c: 3, o: \"29+22=\", n: \"78+13=\" },\n  { r: 19, c: 4, o: \"79+16=\", n: \"51-36=\" }\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\n\nfor (const edit of cellEdits) {\n  const cell = table.getCell(edit.r, edit.c);\n\n  const results = cell.body.search(edit.o, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    // Idempotent re-run guard: if the cell already holds the new value,\n    // there is nothing left to do for it.\n    cell.body.load(\"text\");\n    await context.sync();\n    if (cell.body.text.trim() === edit.n) {\n      continue;\n    }\n    throw new Error(\n      \"Could not find expected text '\" + edit.o + \"' in cell (\" + edit.r + \",\" + edit.c + \").\"\n    );\n  }\n\n  // Replacing the matched range (instead of clearing/re-inserting the whole\n  // cell body) keeps the original run's rFonts/sz/jc formatting intact.\n  results.items[0].insertText(edit.n, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the arithmetic expression text in each cell of the single table\n# in the document body, cell-by-cell (row-major order), preserving all\n# existing run/paragraph formatting. A blanket Find/Replace across the whole\n# document is not safe here because a few expressions (e.g. \"90-53=\", \"33-5=\")\n# repeat in multiple cells but map to different replacement text, so each cell\n# is addressed individually by its 1-based (row, col) position in the table,\n# the existing text is verified, and only the cell's own Range.Text is updated\n# (shrunk by one character to exclude the end-of-cell marker), which keeps the\n# run's existing formatting (font / size / alignment) intact.\n\n$cellEdits = @(\n    @{ r = 1; c = 1; o = '66+18='; n = '26+46=' },\n    @{ r = 1; c = 2; o = '48+38='; n = '39+17=' },\n    @{ r = 1; c = 3; o = '9+48='; n = '85+9=' },\n    @{ r = 1; c = 4; o = '3+78='; n = '19+19=' },\n    @{ r = 1; c = 5; o = '36-18='; n = '3+9=' },\n    @{ r = 2; c = 1; o = '24+58='; n = '50-12=' },\n    @{ r = 2; c = 2; o = '27+27='; n = '97-19=' },\n    @{ r = 2; c = 3; o = '82-74='; n = '55-6=' },\n    @{ r = 2; c = 4; o = '40-35='; n = '61-39=' },\n    @{ r = 2; c = 5; o = '49+44='; n = '29+47=' },\n    @{ r = 3; c = 1; o = '41-17='; n = '66+16=' },\n    @{ r = 3; c = 2; o = '60-59='; n = '92-76=' },\n    @{ r = 3; c = 3; o = '86-67='; n = '38+44=' },\n    @{ r = 3; c = 4; o = '14+18='; n = '73-29=' },\n    @{ r = 3; c = 5; o = '9+37='; n = '91-58=' },\n    @{ r = 4; c = 1; o = '5+7='; n = '40-22=' },\n    @{ r = 4; c = 2; o = '55-47='; n = '93-78=' },\n    @{ r = 4; c = 3; o = '25+38='; n = '85-37=' },\n    @{ r = 4; c = 4; o = '39+43='; n = '59+4=' },\n    @{ r = 4; c = 5; o = '66+19='; n = '9+62=' },\n    @{ r = 5; c = 1; o = '22-15='; n = '71-42=' },\n    @{ r = 5; c = 2; o = '92-18='; n = '57+24=' },\n    @{ r = 5; c = 3; o = '18+57='; n = '19+79=' },\n    @{ r = 5; c = 4; o = '80-29='; n = '19+25=' },\n    @{ r = 5; c = 5; o = '50-28='; n = '4+39=' },\n    @{ r = 6; c = 1; o = '77-38='; n = '78+13=' },\n    @{ r = 6; c = 2; o = '7+77='; n = '48+43=' },\n    @{ r = 6; c = 3; o = '9+64='; n = '45+39=' },\n    @{ r = 6; c = 4; o = '50-48='; n = '17+18=' },\n    @{ r = 6; c = 5; o = '49+35='; n = '16+75=' },\n    @{ r = 7; c = 1; o = '29+45='; n = '61-53=' },\n    @{ r = 7; c = 2; o = '9+18='; n = '63-24=' },\n    @{ r = 7; c = 3; o = '9+58='; n = '19+27=' },\n    @{ r = 7; c = 4; o = '17+55='; n = '43-36=' },\n    @{ r = 7; c = 5; o = '90-53='; n = '81-2=' },\n    @{ r = 8; c = 1; o = '19+66='; n = '36+7=' },\n    @{ r = 8; c = 2; o = '56+37='; n = '47+37=' },\n    @{ r = 8; c = 3; o = '73-36='; n = '26+6=' },\n    @{ r = 8; c = 4; o = '22-9='; n = '77+8=' },\n    @{ r = 8; c = 5; o = '60-8='; n = '26+29=' },\n    @{ r = 9; c = 1; o = '24+67='; n = '94-76=' },\n    @{ r = 9; c = 2; o = '33-5='; n = '9+76=' },\n    @{ r = 9; c = 3; o = '18+38='; n = '72-65=' },\n    @{ r = 9; c = 4; o = '86-8='; n = '33-29=' },\n    @{ r = 9; c = 5; o = '85-67='; n = '83-5=' },\n    @{ r = 10; c = 1; o = '21-7='; n = '84-59=' },\n    @{ r = 10; c = 2; o = '84-65='; n = '37-29=' },\n    @{ r = 10; c = 3; o = '89+9='; n = '73-38=' },\n    @{ r = 10; c = 4; o = '18+65='; n = '95-78=' },\n    @{ r = 10; c = 5; o = '94-86='; n = '74-48=' },\n    @{ r = 11; c = 1; o = '95-79='; n = '46+37=' },\n    @{ r = 11; c = 2; o = '85-36='; n = '90-32=' },\n    @{ r = 11; c = 3; o = '29+26='; n = '74-25=' },\n    @{ r = 11; c = 4; o = '62-38='; n = '6+49=' },\n    @{ r = 11; c = 5; o = '80-7='; n = '48+23=' },\n    @{ r = 12; c = 1; o = '8+13='; n = '61-6=' },\n    @{ r = 12; c = 2; o = '46+9='; n = '69+4=' },\n    @{ r = 12; c = 3; o = '39+42='; n = '61-44=' },\n    @{ r = 12; c = 4; o = '90-14='; n = '14-8=' },\n    @{ r = 12; c = 5; o = '29+35='; n = '82-17=' },\n    @{ r = 13; c = 1; o = '81-6='; n = '30-26=' },\n    @{ r = 13; c = 2; o = '73-9='; n = '9+68=' },\n    @{ r = 13; c = 3; o = '9+5='; n = '62-16=' },\n    @{ r = 13; c = 4; o = '17+49='; n = '9+57=' },\n    @{ r = 13; c = 5; o = '5+48='; n = '25+47=' },\n    @{ r = 14; c = 1; o = '92-3='; n = '19+55=' },\n    @{ r = 14; c = 2; o = '7+7='; n = '58+19=' },\n    @{ r = 14; c = 3; o = '16+25='; n = '70-41=' },\n    @{ r = 14; c = 4; o = '15+79='; n = '86-58=' },\n    @{ r = 14; c = 5; o = '75-7='; n = '59+22=' },\n    @{ r = 15; c = 1; o = '79+13='; n = '81-27=' },\n    @{ r = 15; c = 2; o = '90-53='; n = '36+28=' },\n    @{ r = 15; c = 3; o = '65+16='; n = '45+37=' },\n    @{ r = 15; c = 4; o = '33-6='; n = '47+6=' },\n    @{ r = 15; c = 5; o = '64+9='; n = '90-66=' },\n    @{ r = 16; c = 1; o = '41-37='; n = '33-9=' },\n    @{ r = 16; c = 2; o = '26+36='; n = '41-3=' },\n    @{ r = 16; c = 3; o = '7+87='; n = '57+28=' },\n    @{ r = 16; c = 4; o = '71-17='; n = '51-22=' },\n    @{ r = 16; c = 5; o = '46-28='; n = '16+29=' },\n    @{ r = 17; c = 1; o = '75-48='; n = '68+28=' },\n    @{ r = 17; c = 2; o = '36-7='; n = '23+8=' },\n    @{ r = 17; c = 3; o = '18+26='; n = '88+5=' },\n    @{ r = 17; c = 4; o = '98-29='; n = '80-2=' },\n    @{ r = 17; c = 5; o = '33-5='; n = '11-2=' },\n    @{ r = 18; c = 1; o = '27+24='; n = '64-47=' },\n    @{ r = 18; c = 2; o = '51-36='; n = '16+27=' },\n    @{ r = 18; c = 3; o = '55-27='; n = '70-61=' },\n    @{ r = 18; c = 4; o = '66+17='; n = '84-16=' },\n    @{ r = 18; c = 5; o = '35+47='; n = '87+7=' },\n    @{ r = 19; c = 1; o = '4+87='; n = '51-44=' },\n    @{ r = 19; c = 2; o = '67-29='; n = '38+18=' },\n    @{ r = 19; c = 3; o = '7+39='; n = '40-11=' },\n    @{ r = 19; c = 4; o = '35+26='; n = '61-22=' },\n    @{ r = 19; c = 5; o = '80-57='; n = '71-59=' },\n    @{ r = 20; c = 1; o = '63-56='; n = '18+19=' },\n    @{ r = 20; c = 2; o = '36+25='; n = '37+5=' },\n    @{ r = 20; c = 3; o = '3+28='; n = '30-7=' },\n    @{ r = 20; c = 4; o = '29+22='; n = '78+13=' },\n    @{ r = 20; c = 5; o = '79+16='; n = '51-36=' }\n)\n\n$d = $word.ActiveDocument\n\nif ($d.Tables.Count -lt 1) {\n    throw \"Expected a table in the document body, found none.\"\n}\n\n$tbl = $d.Tables.Item(1)\n\nforeach ($edit in $cellEdits) {\n    $cell = $tbl.Cell($edit.r, $edit.c)\n    $rng = $cell.Range\n    # Drop the trailing cell-mark / paragraph-mark character so only the\n    # visible text is inspected/replaced.\n    $rng.End = $rng.End - 1\n\n    if ($rng.Text -ne $edit.o) {\n        if ($rng.Text -eq $edit.n) {\n            # Already applied (idempotent re-run) - nothing to do.\n            continue\n        }\n        throw \"Unexpected text in cell ($($edit.r),$($edit.c)): expected [$($edit.o)] but found [$($rng.Text)].\"\n    }\n\n    $rng.Text = $edit.n\n}\n"}
